$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Response")

# Insert a new column before the old "productCode" column (B) to hold the
# new StatusCode field; this shifts every existing column one to the right.
$ws.Range("B1").EntireColumn.Insert()

# Match the width of column A (18.6640625 chars) for the newly inserted column.
$ws.Columns.Item(2).ColumnWidth = 17.83

# Header + values for the new StatusCode column.
$ws.Range("B1").Value = "StatusCode"
$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 200

# Restore the view to what was saved with the edit.
$ws.Range("F23").Select()
